# Generate Report for Archive
# - Update status text "Ready for handoff" -> "In Translation" on all
#   sheets that surface per-locale handoff/translation status.
# - Re-fit the status columns now that the label is shorter.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: per-locale status columns are E (zh-cn) and F (de-de).
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

# Per-locale detail sheets: "Status" lives in column C.
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# Columns shrink now that "In Translation" is shorter than
# "Ready for handoff".
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
